$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set values: B1 and A2 are numeric zeros with special styling,
# B2 is the text label "disconnected_elements"
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Apply styling to B1 and A2: bold font, thin box border, centered horizontally, top vertically
foreach ($addr in @("B1", "A2")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
    $cell.Borders.Weight = 2
}
